$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The numeric-looking values in C:F are stored as text (numberStoredAsText),
# so force Text format on each edited cell before assigning, to preserve
# the original "string" cell type instead of converting to a real number.

function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Row 2: runs, balls, fours  (sixes unchanged)
Set-TextValue "C2" "1"
Set-TextValue "D2" "3"
Set-TextValue "E2" "0"

# Row 3: runs, balls, sixes  (fours unchanged)
Set-TextValue "C3" "10"
Set-TextValue "D3" "5"
Set-TextValue "F3" "1"

# Row 4: runs, balls, fours, sixes
Set-TextValue "C4" "4"
Set-TextValue "D4" "2"
Set-TextValue "E4" "1"
Set-TextValue "F4" "0"

# Row 5: runs, balls  (fours, sixes unchanged)
Set-TextValue "C5" "0"
Set-TextValue "D5" "1"
